# Update "想去人数" (want-to-go count) figures for the latest generated output.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 8536
$wsExpo.Range("F4").Value = 6251
$wsExpo.Range("F6").Value = 118
$wsExpo.Range("F9").Value = 338
$wsExpo.Range("F10").Value = 1258

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 8536
$wsAll.Range("F4").Value = 6251
$wsAll.Range("F6").Value = 118
$wsAll.Range("F9").Value = 338
$wsAll.Range("F14").Value = 1259
